$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: NEELU MALIK paid amount changes
$ws.Range("C2").Value = 3662

# Row 3: RUKHSAR KHAN -> VIJAY KHANNA, amount, and LESS THAN DEMAND PAYOUT% 1% -> 0%
$ws.Range("B3").Value = "VIJAY KHANNA"
$ws.Range("C3").Value = 17910
$origStyleE3 = $ws.Range("D3").Style
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0%"
$ws.Range("E3").Style = $origStyleE3

# Row 4: LESS THAN DEMAND -> MORE THAN DEMAND, VIJAY KHANNA -> RUKHSAR KHAN, amount,
# and LESS THAN DEMAND PAYOUT% 2% -> 0%
$ws.Range("A4").Value = "MORE THAN DEMAND"
$ws.Range("B4").Value = "RUKHSAR KHAN"
$ws.Range("C4").Value = 16000
$origStyleE4 = $ws.Range("D4").Style
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0%"
$ws.Range("E4").Style = $origStyleE4

# Rows 5-7 (NEELU MALIK / RUKHSAR KHAN / VIJAY KHANNA "MORE THAN DEMAND" entries)
# are no longer present - delete them entirely so the sheet shrinks to A1:E4
$ws.Range("A5:E7").Delete()
